$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.696.53"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.597.70"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.66"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.70"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "1.822.70"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "1.606.80"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.96"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "26.691.41"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.11"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.00"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.56"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.33"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  +17.90%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.272.34"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.597"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.43"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.17"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.62"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "1.734.64"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.37"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.103"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0512"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.50"
$ws.Range("E51").Value = "  +1.48%  "
